$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data cells to reflect the latest scrape.
# NumberFormat is forced to text ("@") before assignment so that values such as
# "1.00", "0.0000163" etc. are preserved verbatim as strings (matching the
# original inline-string cells) instead of being auto-converted to numbers by
# Excel. The style is then reset to "Normal" so no new cell style is introduced.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.898.82'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -4.53%  '
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.455.39'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -5.98%  '
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '545.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -5.12%  '
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -7.21%  '
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.605'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.08%  '
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.453.40'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.97%  '
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.107'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -10.05%  '
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.96%  '
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.34'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -8.45%  '
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.353'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -7.39%  '
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -8.14%  '
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.895.76'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -6.14%  '
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000163'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -9.72%  '
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.784.62'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -4.55%  '
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.470.56'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -6.03%  '
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.03'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -8.22%  '
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -8.55%  '
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.16'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -8.11%  '
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '317.60'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -7.53%  '
$ws.Range("E22").Style = "Normal"

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -6.61%  '
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.19%  '
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'PEPE'
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0₃0972'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -10.80%  '
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.581.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -5.70%  '
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.50'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.93%  '
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '535.57'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -9.99%  '
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -9.55%  '
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.64'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.34%  '
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.148'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -8.37%  '
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.89'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -8.60%  '
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -10.03%  '
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.87'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -10.98%  '
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.87'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -9.25%  '
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = "Normal"

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -6.95%  '
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.29'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -7.30%  '
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '143.95'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.66%  '
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.70'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -9.35%  '
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.98'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.56%  '
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.29'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -9.52%  '
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '146.35'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -7.02%  '
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -8.80%  '
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.82'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -12.46%  '
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0529'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -10.32%  '
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0939'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -6.13%  '
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.581'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -7.89%  '
$ws.Range("E51").Style = "Normal"
